# DOMA-6936: add "Верифицирован" (Is verified) column (H) to the contacts
# import example sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Copy the formatting (style) of column G into the new column H ---
# This gives H1 the header style (s=2), H2/H3/H5/H6 the data style (s=2)
# and H4 the empty-but-styled cell (s=2), plus H7:H50 the blank-row style
# (s=6) — exactly mirroring column G's per-row styles without touching any
# values.
$ws.Range("G1:G50").Copy()
$ws.Range("H1:H50").PasteSpecial(-4122)

# --- 2. Match column H's width to columns F:G (17.6719) as closely as the
# column-width API allows. ---
$ws.Columns.Item(8).ColumnWidth = 16.75

# --- 3. Fill in the new column's header + data values ---
$ws.Range("H1").Value = "Верифицирован"
$ws.Range("H2").Value = "Да"
$ws.Range("H3").Value = "Нет"
$ws.Range("H4").Value = ""
$ws.Range("H5").Value = "да"
$ws.Range("H6").Value = "нет"

# --- 4. Update the F5 hyperlink's displayed text to match the cell's real
# text ("ttest@example.com" instead of the stale "test@example.com").
# The host does not support editing a Hyperlink in place, so rebuild the
# full set (targets/display text unchanged except for F5). ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:test@example.com", "", "", "test@example.com")
$ws.Hyperlinks.Add($ws.Range("F3"), "mailto:test@example.com", "", "", "test1@example.com")
$ws.Hyperlinks.Add($ws.Range("F4"), "mailto:test@example.com", "", "", "test2@example.com")
$ws.Hyperlinks.Add($ws.Range("F5"), "mailto:test@example.com", "", "", "ttest@example.com")

# Adding a hyperlink re-styles the anchor cell with the built-in "Hyperlink"
# style (underline + theme color). Restore the original per-cell style
# (s=2, same as before the edit) by pasting formats back from the
# untouched G column.
$ws.Range("G2").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("G3").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("G4").Copy()
$ws.Range("F4").PasteSpecial(-4122)
$ws.Range("G5").Copy()
$ws.Range("F5").PasteSpecial(-4122)
